# This script applies the "Update countries & provincias Spain" data refresh
# to the Pais worksheet: it updates the last-updated timestamp and refreshes
# the per-country statistics (columns B-H), including three countries whose
# case counts grew enough to move them up one rank (Paises Bajos, Islandia,
# Liechtenstein), which pushes the country that used to occupy that row down
# by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Octubre de 2020 a las 14:40"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 8038391
$ws.Range("C4").Value = 602
$ws.Range("D4").Value = 5185991
$ws.Range("E4").Value = 2632379
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 220021

# Row 29: Paises Bajos
$ws.Range("A29").Value = "Paises Bajos"
$ws.Range("B29").Value = 188876
$ws.Range("C29").Value = 7378
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 35
$ws.Range("H29").Value = 6631

# Row 30: Canada
$ws.Range("A30").Value = "Canada"
$ws.Range("B30").Value = 182839
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 154258
$ws.Range("E30").Value = 18954
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 9627

# Row 42: Kuwait
$ws.Range("A42").Value = "Kuwait"
$ws.Range("B42").Value = 112737
$ws.Range("C42").Value = 844
$ws.Range("D42").Value = 104508
$ws.Range("E42").Value = 7557
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 8
$ws.Range("H42").Value = 672

# Row 47: Suecia
$ws.Range("A47").Value = "Suecia"
$ws.Range("B47").Value = 100654
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 6
$ws.Range("H47").Value = 5899

# Row 73: Azerbaiyan
$ws.Range("A73").Value = "Azerbaiyan"
$ws.Range("B73").Value = 42381
$ws.Range("C73").Value = 277
$ws.Range("D73").Value = 39468
$ws.Range("E73").Value = 2301
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 612

# Row 78: Dinamarca
$ws.Range("A78").Value = "Dinamarca"
$ws.Range("B78").Value = 33101
$ws.Range("C78").Value = 290
$ws.Range("D78").Value = 27225
$ws.Range("E78").Value = 5202
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 674

# Row 89: Republica de Macedonia
$ws.Range("A89").Value = "Republica de Macedonia"
$ws.Range("B89").Value = 21193
$ws.Range("C89").Value = 80
$ws.Range("D89").Value = 16397
$ws.Range("E89").Value = 3996
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 800

# Row 95: Noruega
$ws.Range("A95").Value = "Noruega"
$ws.Range("B95").Value = 15639
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 11863
$ws.Range("E95").Value = 3499
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 277

# Row 127: Hong Kong
$ws.Range("A127").Value = "Hong Kong"
$ws.Range("B127").Value = 5202
$ws.Range("C127").Value = 8
$ws.Range("D127").Value = 4931
$ws.Range("E127").Value = 166
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 105

# Row 143: Islandia
$ws.Range("A143").Value = "Islandia"
$ws.Range("B143").Value = 3668
$ws.Range("C143").Value = 86
$ws.Range("D143").Value = 2587
$ws.Range("E143").Value = 1071
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 10

# Row 144: Tailandia
$ws.Range("A144").Value = "Tailandia"
$ws.Range("B144").Value = 3643
$ws.Range("C144").Value = 2
$ws.Range("D144").Value = 3457
$ws.Range("E144").Value = 127
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 59

# Row 145: Gambia
$ws.Range("A145").Value = "Gambia"
$ws.Range("B145").Value = 3636
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 2593
$ws.Range("E145").Value = 925
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 118

# Row 168: Vietnam
$ws.Range("A168").Value = "Vietnam"
$ws.Range("B168").Value = 1113
$ws.Range("C168").Value = 3
$ws.Range("D168").Value = 1025
$ws.Range("E168").Value = 53
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 35

# Row 181: Islas Feroe
$ws.Range("A181").Value = "Islas Feroe"
$ws.Range("B181").Value = 477
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 465
$ws.Range("E181").Value = 12
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0

# Row 193: Liechtenstein
$ws.Range("A193").Value = "Liechtenstein"
$ws.Range("B193").Value = 163
$ws.Range("C193").Value = 15
$ws.Range("D193").Value = 128
$ws.Range("E193").Value = 34
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 1

# Row 194: Bonaire, San Eustaquio y Saba
$ws.Range("A194").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("B194").Value = 148
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 78
$ws.Range("E194").Value = 68
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 2
